$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename "2022 ..." row labels to "2020 ..." ---
# (Texas grid-cost refresh brought the cost-basis year label back to 2020)
$ws.Range("A9").Value  = "2020 CapEx"
$ws.Range("A15").Value = "2020 OpEx ($/kw-yr)"
$ws.Range("A22").Value = "2020 PV base installed cost"
$ws.Range("A28").Value = "2020 PV OpEx"

# --- Grid interconnection cost updates (rows 38-40) ---
$ws.Range("B38").Value = 20.91
$ws.Range("C38").Value = 25.87
$ws.Range("D38").Value = 47.26
$ws.Range("E38").Value = 30.91
$ws.Range("F38").Value = 49.34

$ws.Range("B39").Value = 33.34
$ws.Range("C39").Value = 79.68
$ws.Range("D39").Value = 61.68
$ws.Range("E39").Value = 23.72
$ws.Range("F39").Value = 67.43

$ws.Range("B40").Value = 24.04
$ws.Range("C40").Value = 95.83
$ws.Range("D40").Value = 52.38
$ws.Range("E40").Value = 65.7
$ws.Range("F40").Value = 0.68

# --- Updated Texas site (Site 2, column C) location ---
# Re-use the existing B38 cell format (center-aligned "General") for C4
# instead of its old font/format, matching the style Excel picked here.
$ws.Range("B38").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "32.318714, -100.18"
$ws.Range("C4").Value = 32.318714
$ws.Range("C5").Value = -100.18

# --- Reflect the last-selected cell (C4, the edited coordinate) ---
$ws.Range("C4").Select() | Out-Null
